$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58 - this shifts the existing rows 58:90
# down to 59:91 (and carries the D-column date style/format along, since
# Excel's Insert() copies formatting from the row above by default).
$ws.Rows.Item(58).Insert()

# Populate the freshly inserted row 58 with the new weekly record.
$ws.Cells.Item(58, 1).Value = 6
$ws.Cells.Item(58, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(58, 3).Value = "Metropolitana"
$ws.Cells.Item(58, 4).Value = 44438
$ws.Cells.Item(58, 5).Value = 13
$ws.Cells.Item(58, 6).Value = 100112001
$ws.Cells.Item(58, 7).Value = "Berenjena"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 350
$ws.Cells.Item(58, 11).Value = 7000
$ws.Cells.Item(58, 12).Value = 8000
$ws.Cells.Item(58, 13).Value = 7429
$ws.Cells.Item(58, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(58, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(58, 16).Value = 149
$ws.Cells.Item(58, 17).Value = 50
$ws.Cells.Item(58, 18).Value = "Hortaliza"
